# Reorders classFields rows to reflect the new field declaration order
# (standard relationship between microservices / MSM measure work).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

$ws.Cells.Item(2, 2).Value = "`$assertionsDisabled"
$ws.Cells.Item(2, 4).Value = "boolean"
$ws.Cells.Item(3, 2).Value = "serialVersionUID"
$ws.Cells.Item(3, 4).Value = "long"
$ws.Cells.Item(5, 2).Value = "`$assertionsDisabled"
$ws.Cells.Item(5, 4).Value = "boolean"
$ws.Cells.Item(6, 2).Value = "serialVersionUID"
$ws.Cells.Item(6, 4).Value = "long"
$ws.Cells.Item(8, 2).Value = "LOW"
$ws.Cells.Item(10, 2).Value = "HIGH"
$ws.Cells.Item(11, 2).Value = "HIGHEST"
$ws.Cells.Item(13, 2).Value = "MEDIUM"
$ws.Cells.Item(15, 2).Value = "isConfirmed"
$ws.Cells.Item(15, 4).Value = "java.lang.Boolean"
$ws.Cells.Item(16, 2).Value = "subscriptionDate"
$ws.Cells.Item(16, 4).Value = "java.time.LocalDateTime"
$ws.Cells.Item(17, 2).Value = "`$assertionsDisabled"
$ws.Cells.Item(17, 4).Value = "boolean"
$ws.Cells.Item(18, 2).Value = "serialVersionUID"
$ws.Cells.Item(18, 4).Value = "long"
$ws.Cells.Item(20, 2).Value = "observedUsers"
$ws.Cells.Item(20, 4).Value = "java.util.Set"
$ws.Cells.Item(21, 2).Value = "observingUsers"
$ws.Cells.Item(21, 4).Value = "java.util.Set"
$ws.Cells.Item(22, 2).Value = "imageUrl"
$ws.Cells.Item(22, 4).Value = "java.lang.String"
$ws.Cells.Item(23, 2).Value = "key"
$ws.Cells.Item(24, 2).Value = "username"
$ws.Cells.Item(24, 4).Value = "java.lang.String"
$ws.Cells.Item(25, 2).Value = "pageSize"
$ws.Cells.Item(25, 4).Value = "java.lang.Integer"
$ws.Cells.Item(26, 2).Value = "query"
$ws.Cells.Item(26, 4).Value = "java.lang.String"
$ws.Cells.Item(27, 2).Value = "pageNumber"
$ws.Cells.Item(28, 2).Value = "`$assertionsDisabled"
$ws.Cells.Item(28, 4).Value = "boolean"
$ws.Cells.Item(29, 2).Value = "serialVersionUID"
$ws.Cells.Item(29, 4).Value = "long"
$ws.Cells.Item(33, 2).Value = "phoneNumber"
$ws.Cells.Item(34, 2).Value = "street"
$ws.Cells.Item(35, 2).Value = "profileImageUrl"
$ws.Cells.Item(36, 2).Value = "dateOfBirth"
$ws.Cells.Item(37, 2).Value = "description"
$ws.Cells.Item(38, 2).Value = "city"
$ws.Cells.Item(39, 2).Value = "gender"
$ws.Cells.Item(40, 2).Value = "firstName"
$ws.Cells.Item(41, 2).Value = "lastName"
$ws.Cells.Item(42, 2).Value = "country"
$ws.Cells.Item(43, 2).Value = "postalCode"
$ws.Cells.Item(44, 2).Value = "query"
$ws.Cells.Item(44, 4).Value = "java.lang.String"
$ws.Cells.Item(45, 2).Value = "pageNumber"
$ws.Cells.Item(46, 2).Value = "pageSize"
$ws.Cells.Item(46, 4).Value = "java.lang.Integer"
$ws.Cells.Item(47, 2).Value = "`$assertionsDisabled"
$ws.Cells.Item(47, 4).Value = "boolean"
$ws.Cells.Item(48, 2).Value = "serialVersionUID"
$ws.Cells.Item(48, 4).Value = "long"
$ws.Cells.Item(50, 2).Value = "FORUM"
$ws.Cells.Item(51, 2).Value = "PRODUCT"
$ws.Cells.Item(53, 2).Value = "ORDER"
$ws.Cells.Item(54, 2).Value = "`$assertionsDisabled"
$ws.Cells.Item(54, 4).Value = "boolean"
$ws.Cells.Item(55, 2).Value = "serialVersionUID"
$ws.Cells.Item(55, 4).Value = "long"
$ws.Cells.Item(57, 2).Value = "acknowledgedUsers"
$ws.Cells.Item(57, 4).Value = "java.util.Set"
$ws.Cells.Item(58, 2).Value = "affectedUsers"
$ws.Cells.Item(58, 4).Value = "java.util.Set"
$ws.Cells.Item(59, 2).Value = "relatedId"
$ws.Cells.Item(60, 2).Value = "priority"
$ws.Cells.Item(60, 4).Value = "org.andante.activity.enums.Priority"
$ws.Cells.Item(61, 2).Value = "description"
$ws.Cells.Item(62, 2).Value = "eventTimestamp"
$ws.Cells.Item(62, 4).Value = "java.time.LocalDateTime"
$ws.Cells.Item(63, 2).Value = "id"
$ws.Cells.Item(63, 4).Value = "java.lang.String"
$ws.Cells.Item(64, 2).Value = "domain"
$ws.Cells.Item(64, 4).Value = "org.andante.activity.enums.Domain"
$ws.Cells.Item(65, 2).Value = "`$assertionsDisabled"
$ws.Cells.Item(65, 4).Value = "boolean"
$ws.Cells.Item(66, 2).Value = "serialVersionUID"
$ws.Cells.Item(66, 4).Value = "long"
$ws.Cells.Item(68, 2).Value = "affectedUsers"
$ws.Cells.Item(68, 4).Value = "java.util.Set"
$ws.Cells.Item(70, 2).Value = "relatedId"
$ws.Cells.Item(70, 4).Value = "java.lang.String"
$ws.Cells.Item(71, 2).Value = "acknowledgedUsers"
$ws.Cells.Item(71, 4).Value = "java.util.Set"
$ws.Cells.Item(72, 2).Value = "description"
$ws.Cells.Item(72, 4).Value = "java.lang.String"
$ws.Cells.Item(73, 2).Value = "priority"
$ws.Cells.Item(73, 4).Value = "org.andante.activity.enums.Priority"
$ws.Cells.Item(74, 2).Value = "id"
$ws.Cells.Item(75, 2).Value = "domain"
$ws.Cells.Item(75, 4).Value = "org.andante.activity.enums.Domain"
$ws.Cells.Item(76, 2).Value = "imageUrl"
$ws.Cells.Item(77, 2).Value = "username"
$ws.Cells.Item(78, 2).Value = "subscriptionDate"
$ws.Cells.Item(78, 4).Value = "java.time.LocalDateTime"
$ws.Cells.Item(80, 2).Value = "emailAddress"
$ws.Cells.Item(80, 4).Value = "java.lang.String"
$ws.Cells.Item(81, 2).Value = "profileImageUrl"
$ws.Cells.Item(82, 2).Value = "city"
$ws.Cells.Item(83, 2).Value = "lastName"
$ws.Cells.Item(84, 2).Value = "street"
$ws.Cells.Item(85, 2).Value = "gender"
$ws.Cells.Item(86, 2).Value = "phoneNumber"
$ws.Cells.Item(87, 2).Value = "postalCode"
$ws.Cells.Item(88, 2).Value = "description"
$ws.Cells.Item(89, 2).Value = "firstName"
$ws.Cells.Item(90, 2).Value = "dateOfBirth"
$ws.Cells.Item(91, 2).Value = "country"
$ws.Cells.Item(92, 2).Value = "key"
$ws.Cells.Item(94, 2).Value = "username"
$ws.Cells.Item(97, 2).Value = "`$assertionsDisabled"
$ws.Cells.Item(97, 4).Value = "boolean"
$ws.Cells.Item(98, 2).Value = "serialVersionUID"
$ws.Cells.Item(98, 4).Value = "long"
